# Auto-generated edit script for QuizTestAppData.xlsx
# Applies: row171 answer update, 30 new QuestionBank rows (Su Dia 7),
# Category sheet counters, and Quizzes sheet question-list update.

$wb = $excel.ActiveWorkbook
$wsQB = $wb.Worksheets.Item("QuestionBank")
$wsCat = $wb.Worksheets.Item("Category")
$wsQuiz = $wb.Worksheets.Item("Quizzes")

# --- Update existing row 171 ("Where is bug") with new answer data ---
$wsQB.Cells.Item(171,3).Value = "A, B, C"
$wsQB.Cells.Item(171,4).Value = "A. In here`nB. In some where`nC. Ở đâu đó"
$wsQB.Cells.Item(171,5).Value = "N`nN`nN`nN`nN`nN`nN`nN`nN`nN`n"

# --- Append 30 new question rows (204-233), category "Su Dia 7" ---
$wsQB.Cells.Item(204,1).Value = "Thành phố nào là thủ đô của Brazil?"
$wsQB.Cells.Item(204,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(204,3).Value = "A"
$wsQB.Cells.Item(204,4).Value = "A. Brasília`nB. Rio de Janeiro`nC. São Paulo"
$wsQB.Cells.Item(205,1).Value = "Dòng sông nổi tiếng chảy qua thung lũng lịch sử nổi tiếng Mesopotamia là:"
$wsQB.Cells.Item(205,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(205,3).Value = "C"
$wsQB.Cells.Item(205,4).Value = "A. Sông Amazon`nB. Sông Nile`nC. Sông Eufrat"
$wsQB.Cells.Item(206,1).Value = "Đồng bằng nổi tiếng và phù hợp cho nông nghiệp nằm ở miền Bắc Việt Nam là:"
$wsQB.Cells.Item(206,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(206,3).Value = "B"
$wsQB.Cells.Item(206,4).Value = "A. Đồng bằng Sông Cửu Long`nB. Đồng bằng Sông Hồng`nC. Đồng bằng Sông Mê Kông"
$wsQB.Cells.Item(207,1).Value = "Quốc gia nào không nằm ở châu Á?"
$wsQB.Cells.Item(207,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(207,3).Value = "C"
$wsQB.Cells.Item(207,4).Value = "A. Nhật Bản`nB. Ấn Độ`nC. Brazil"
$wsQB.Cells.Item(208,1).Value = "Thành phố lớn nhất nước Mỹ là:"
$wsQB.Cells.Item(208,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(208,3).Value = "B"
$wsQB.Cells.Item(208,4).Value = "A. New York`nB. Los Angeles`nC. Chicago"
$wsQB.Cells.Item(209,1).Value = "Quốc gia nằm ở Đông Nam Á, có thủ đô là Jakarta là:"
$wsQB.Cells.Item(209,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(209,3).Value = "A"
$wsQB.Cells.Item(209,4).Value = "A. Indonesia`nB. Malaysia`nC. Philippines"
$wsQB.Cells.Item(210,1).Value = "Đỉnh núi nổi tiếng ở châu Phi là:"
$wsQB.Cells.Item(210,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(210,3).Value = "B"
$wsQB.Cells.Item(210,4).Value = "A. Everest`nB. Kilimanjaro`nC. Mont Blanc"
$wsQB.Cells.Item(211,1).Value = "Quốc gia nào nằm ở Bán đảo Iberia?"
$wsQB.Cells.Item(211,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(211,3).Value = "A"
$wsQB.Cells.Item(211,4).Value = "A. Tây Ban Nha`nB. Italia`nC. Hy Lạp"
$wsQB.Cells.Item(212,1).Value = "Đảo nổi tiếng với tượng chúa Kitô Vua là:"
$wsQB.Cells.Item(212,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(212,3).Value = "C"
$wsQB.Cells.Item(212,4).Value = "A. Đảo Hawaii`nB. Đảo Jeju`nC. Đảo Corcovado"
$wsQB.Cells.Item(213,1).Value = "Thành phố nổi tiếng với Cung điện Buckingham là:"
$wsQB.Cells.Item(213,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(213,3).Value = "B"
$wsQB.Cells.Item(213,4).Value = "A. Paris`nB. London`nC. Moscow"
$wsQB.Cells.Item(214,1).Value = "Quốc gia nào có dân số đông nhất thế giới?"
$wsQB.Cells.Item(214,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(214,3).Value = "A"
$wsQB.Cells.Item(214,4).Value = "A. Trung Quốc`nB. Ấn Độ`nC. Mỹ"
$wsQB.Cells.Item(215,1).Value = "Dòng sông Mississippi chảy qua bao nhiêu bang của Mỹ?"
$wsQB.Cells.Item(215,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(215,3).Value = "B"
$wsQB.Cells.Item(215,4).Value = "A. 1`nB. 2`nC. 3"
$wsQB.Cells.Item(216,1).Value = "Thành phố nào là thủ đô của Ý?"
$wsQB.Cells.Item(216,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(216,3).Value = "A"
$wsQB.Cells.Item(216,4).Value = "A. Roma`nB. Milano`nC. Napoli"
$wsQB.Cells.Item(217,1).Value = "Đảo nằm giữa biển Địa Trung Hải và được biết đến với cung điện Alhambra là:"
$wsQB.Cells.Item(217,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(217,3).Value = "C"
$wsQB.Cells.Item(217,4).Value = "A. Đảo Balearic`nB. Đảo Corsica`nC. Đảo Sicily"
$wsQB.Cells.Item(218,1).Value = "Hệ thống đường sắt lớn nhất châu Âu là:"
$wsQB.Cells.Item(218,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(218,3).Value = "C"
$wsQB.Cells.Item(218,4).Value = "A. Eurostar`nB. TGV`nC. Trans-Siberian"
$wsQB.Cells.Item(219,1).Value = "Quốc gia nào là quê hương của đền Taj Mahal?"
$wsQB.Cells.Item(219,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(219,3).Value = "A"
$wsQB.Cells.Item(219,4).Value = "A. Ấn Độ`nB. Trung Quốc`nC. Nhật Bản"
$wsQB.Cells.Item(220,1).Value = "Dải đất hẹp nằm giữa biển Đen và biển Địa Trung Hải là:"
$wsQB.Cells.Item(220,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(220,3).Value = "A"
$wsQB.Cells.Item(220,4).Value = "A. Bán đảo Balkan`nB. Bán đảo Iberia`nC. Bán đảo Scandinavia"
$wsQB.Cells.Item(221,1).Value = "Núi lửa Vesuvius nằm ở quốc gia nào?"
$wsQB.Cells.Item(221,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(221,3).Value = "A"
$wsQB.Cells.Item(221,4).Value = "A. Ý`nB. Tây Ban Nha`nC. Hy Lạp"
$wsQB.Cells.Item(222,1).Value = "Quốc gia nào là quê hương của vịnh Hạ Long?"
$wsQB.Cells.Item(222,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(222,3).Value = "A"
$wsQB.Cells.Item(222,4).Value = "A. Việt Nam`nB. Thái Lan`nC. Indonesia"
$wsQB.Cells.Item(223,1).Value = "Dải đất hẹp nằm giữa biển Caribe và Đại Tây Dương là:"
$wsQB.Cells.Item(223,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(223,3).Value = "B"
$wsQB.Cells.Item(223,4).Value = "A. Bán đảo Florida`nB. Bán đảo Yucatan`nC. Bán đảo Labrador"
$wsQB.Cells.Item(224,1).Value = "Thành phố nào được gọi là `"thành phố không bao giờ ngủ`"?"
$wsQB.Cells.Item(224,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(224,3).Value = "B"
$wsQB.Cells.Item(224,4).Value = "A. Tokyo`nB. New York`nC. London"
$wsQB.Cells.Item(225,1).Value = "Đồng bằng lớn nhất thế giới là:"
$wsQB.Cells.Item(225,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(225,3).Value = "C"
$wsQB.Cells.Item(225,4).Value = "A. Đồng bằng Sông Cửu Long`nB. Đồng bằng Ganges-Brahmaputra`nC. Đồng bằng Amazon"
$wsQB.Cells.Item(226,1).Value = "Quốc gia nào là quê hương của vùng đất tuyết Alaska?"
$wsQB.Cells.Item(226,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(226,3).Value = "B"
$wsQB.Cells.Item(226,4).Value = "A. Canada`nB. Mỹ`nC. Nga"
$wsQB.Cells.Item(227,1).Value = "Đảo nổi tiếng với khối núi Phật A Di Đà là:"
$wsQB.Cells.Item(227,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(227,3).Value = "B"
$wsQB.Cells.Item(227,4).Value = "A. Đảo Bali`nB. Đảo Lombok`nC. Đảo Jeju"
$wsQB.Cells.Item(228,1).Value = "Dòng sông nào là dòng sông dài nhất châu Á?"
$wsQB.Cells.Item(228,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(228,3).Value = "C"
$wsQB.Cells.Item(228,4).Value = "A. Sông Mekong`nB. Sông Ganges`nC. Sông Yangtze"
$wsQB.Cells.Item(229,1).Value = "Quốc gia nào nằm ở Bắc Âu?"
$wsQB.Cells.Item(229,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(229,3).Value = "A"
$wsQB.Cells.Item(229,4).Value = "A. Phần Lan`nB. Hà Lan`nC. Đan Mạch"
$wsQB.Cells.Item(230,1).Value = "Quốc gia nào nằm ở khu vực Đông Nam Á?"
$wsQB.Cells.Item(230,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(230,3).Value = "B"
$wsQB.Cells.Item(230,4).Value = "A. Argentina`nB. Malaysia`nC. Ba Lan"
$wsQB.Cells.Item(231,1).Value = "Đại ngàn nổi tiếng của Việt Nam là:"
$wsQB.Cells.Item(231,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(231,3).Value = "B"
$wsQB.Cells.Item(231,4).Value = "A. Đại Sơn`nB. Trường Sơn`nC. Trường Giang"
$wsQB.Cells.Item(232,1).Value = "Nước nào là quốc gia lớn nhất thế giới về diện tích?"
$wsQB.Cells.Item(232,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(232,3).Value = "A"
$wsQB.Cells.Item(232,4).Value = "A. Nga`nB. Trung Quốc`nC. Mỹ"
$wsQB.Cells.Item(233,1).Value = "Kênh đào nổi tiếng nằm ở Ai Cập là:"
$wsQB.Cells.Item(233,2).Value = "root/Course : Junior high school/Sử Địa 7"
$wsQB.Cells.Item(233,3).Value = "B"
$wsQB.Cells.Item(233,4).Value = "A. Kênh Panama`nB. Kênh Suez`nC. Kênh Kiel"

# --- Update Category sheet counters ---
$wsCat.Cells.Item(1,2).Value = 233
$wsCat.Cells.Item(13,2).Value = 55
$wsCat.Cells.Item(15,2).Value = 30

# --- Update Quizzes sheet question list (add row 170 to quiz) ---
# Force text (avoid Excel auto-parsing the comma-separated digits as a number),
# then restore the default "Normal" style so no stray number format remains.
$wsQuiz.Cells.Item(7,5).NumberFormat = "@"
$wsQuiz.Cells.Item(7,5).Value = ",171,202,170"
$wsQuiz.Cells.Item(7,5).Style = "Normal"

